$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'30.239.75"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -0.34%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.863.82"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -0.32%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'0.9990"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.22%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'236.03"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +0.27%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'0.9992"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -0.15%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.4703"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.37%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.2902"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +1.89%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.06567"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +0.39%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'21.84"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +2.07%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  +1.46%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'97.80"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -0.42%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.862.40"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -0.40%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'5.117"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +0.10%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.6796"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +0.45%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'  -3.39%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'30.211.28"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'13.63"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +6.88%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'0.000007649"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +4.66%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.9991"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'2.105.94"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -0.65%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.9989"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -0.27%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'5.235"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -4.50%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'6.197"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +0.69%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'167.30"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +1.14%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'9.202"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.39%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'18.97"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -1.09%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'1.954"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  +1.15%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'1.371"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'0.09936"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +2.43%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'4.345"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -1.30%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'1.467"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.56%  "
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  -1.40%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'0.04713"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'1.124"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.54%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'0.7036"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -0.49%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'2.703"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.88%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.01878"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +0.83%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'2.604"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +2.60%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'6.337"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'73.61"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -1.38%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'1.941"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.54%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.8401"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -1.25%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'103.79"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.9986"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.4152"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -0.96%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'7.068"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -2.18%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'9.155"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -1.35%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'933.93"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -0.31%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'34.13"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  -0.38%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.05661"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  +0.34%  "
$c.Style = "Normal"
